# Updated cryptos list on Mon Jul  3 13:15:34 UTC 2023 with GitHub Actions
#
# Refreshes the live price (col D) / 1h-volume-change (col E) snapshot for
# each coin row on Sheet1, and fixes two coin pairs whose rows had swapped
# rank order (Solana<->Litecoin at rows 10/11, ShibaInu<->Uniswap at rows
# 20/21, RenderToken<->FraxShare at rows 39/40, Quant<->TheSandbox<->
# TrustWalletToken at rows 41/42/43) by rewriting name/link/price/volume
# together for those rows.
#
# Price cells (col D) are free-form text in this sheet (e.g. "30.587.00",
# "109.00", "0.07710") rather than real numbers, so for any value Excel
# would otherwise auto-convert to a Number (dropping trailing zeros /
# re-grouping the digits), we force the cell to Text via NumberFormat="@"
# before writing it, to keep the exact display string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.587.00'
$ws.Range('E2').Value = '  +0.17%  '
$ws.Range('D3').Value = '1.960.96'
$ws.Range('E3').Value = '  +2.26%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9977'
$ws.Range('E4').Value = '  -0.25%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '247.57'
$ws.Range('E5').Value = '  +0.82%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9979'
$ws.Range('E6').Value = '  -0.23%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4833'
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2924'
$ws.Range('E8').Value = '  +0.87%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06754'
$ws.Range('E9').Value = '  +0.47%  '
$ws.Range('B10').Value = 'Litecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '109.00'
$ws.Range('E10').Value = '  -2.22%  '
$ws.Range('B11').Value = 'Solana'
$ws.Range('C11').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '19.23'
$ws.Range('E11').Value = '  +1.46%  '
$ws.Range('D12').Value = '1.964.06'
$ws.Range('E12').Value = '  +2.44%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07710'
$ws.Range('E13').Value = '  +1.96%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.419'
$ws.Range('E14').Value = '  +2.38%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6840'
$ws.Range('E15').Value = '  +1.86%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '289.75'
$ws.Range('E16').Value = '  -3.24%  '
$ws.Range('D17').Value = '30.589.67'
$ws.Range('E17').Value = '  +0.22%  '
$ws.Range('E18').Value = '  +1.06%  '
$ws.Range('D19').Value = '2.220.51'
$ws.Range('E19').Value = '  +2.53%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.629'
$ws.Range('E20').Value = '  +1.01%  '
$ws.Range('B21').Value = 'ShibaInu'
$ws.Range('C21').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.000007659'
$ws.Range('E21').Value = '  +1.01%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9973'
$ws.Range('E22').Value = '  -0.29%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9976'
$ws.Range('E23').Value = '  -0.26%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.586'
$ws.Range('E24').Value = '  +2.40%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.849'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '170.88'
$ws.Range('E26').Value = '  +3.43%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.93'
$ws.Range('E27').Value = '  -1.89%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.173'
$ws.Range('E28').Value = '  +2.81%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1066'
$ws.Range('E29').Value = '  +0.34%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.437'
$ws.Range('E30').Value = '  +0.44%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.752'
$ws.Range('E31').Value = '  +16.17%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.472'
$ws.Range('E32').Value = '  +7.81%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05068'
$ws.Range('E33').Value = '  +1.18%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7674'
$ws.Range('E34').Value = '  +3.72%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.174'
$ws.Range('E35').Value = '  +3.01%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.725'
$ws.Range('E36').Value = '  -0.10%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02025'
$ws.Range('E37').Value = '  +0.07%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.713'
$ws.Range('E38').Value = '  +1.22%  '
$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.419'
$ws.Range('E39').Value = '  +9.88%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.121'
$ws.Range('E40').Value = '  +5.09%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8813'
$ws.Range('E41').Value = '  +2.17%  '
$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '109.35'
$ws.Range('E42').Value = '  -1.14%  '
$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4442'
$ws.Range('E43').Value = '  -0.02%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '69.86'
$ws.Range('E44').Value = '  -1.91%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.9979'
$ws.Range('E45').Value = '  -0.23%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.468'
$ws.Range('E46').Value = '  +3.12%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.1268'
$ws.Range('E47').Value = '  +2.72%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.299'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '35.76'
$ws.Range('E49').Value = '  +2.55%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '924.56'
$ws.Range('E50').Value = '  +6.72%  '
$ws.Range('E51').Value = '  -3.97%  '
